# Trade #9 closed at 2026-02-16 22:52:51 - base_strategy DOWN +0.000%
# Append a new trade row (row 10) to both the "All Trades" and the
# "base_strategy" sheets, mirroring the existing OPEN trade rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3 already has the exact same shape of data we need (same side,
    # prices, status, capital, slippage, confidence and entry reason), so
    # duplicate it into row 10 and then fix up the trade number and time.
    $ws.Range("A3:Q3").Copy($ws.Range("A10:Q10"))

    $ws.Cells.Item(10, 1).Value = 9
    $ws.Cells.Item(10, 3).Value = "22:52:51"
}
